$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "68.237.50"
Set-TextValue "E2" "  +7.38%  "
Set-TextValue "D3" "3.672.01"
Set-TextValue "E3" "  +5.64%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "421.14"
Set-TextValue "E5" "  +1.45%  "
Set-TextValue "D6" "129.93"
Set-TextValue "E6" "  +0.72%  "
Set-TextValue "D7" "0.654"
Set-TextValue "E7" "  +3.31%  "
Set-TextValue "D8" "3.661.05"
Set-TextValue "E8" "  +5.53%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "0.768"
Set-TextValue "E10" "  +2.29%  "
Set-TextValue "D11" "0.198"
Set-TextValue "E11" "  +26.68%  "
Set-TextValue "D12" "0.0000443"
Set-TextValue "E12" "  +94.45%  "
Set-TextValue "D13" "42.13"
Set-TextValue "E13" "  -0.54%  "
Set-TextValue "D14" "9.83"
Set-TextValue "E14" "  +1.55%  "
Set-TextValue "D15" "4.226.92"
Set-TextValue "E15" "  +5.00%  "
Set-TextValue "E16" "  +0.35%  "
Set-TextValue "D17" "3.685.43"
Set-TextValue "E17" "  +5.70%  "
Set-TextValue "D18" "20.08"
Set-TextValue "E18" "  -0.81%  "
Set-TextValue "E19" "  +2.10%  "
Set-TextValue "D20" "68.032.31"
Set-TextValue "E20" "  +7.25%  "
Set-TextValue "D21" "12.48"
Set-TextValue "E21" "  +1.03%  "
Set-TextValue "D22" "460.57"
Set-TextValue "E22" "  +0.42%  "
Set-TextValue "D23" "89.13"
Set-TextValue "E23" "  -0.75%  "
Set-TextValue "D24" "13.57"
Set-TextValue "E24" "  +3.38%  "
Set-TextValue "D25" "3.04"
Set-TextValue "E25" "  -6.94%  "
Set-TextValue "B26" "EthereumClassic"
Set-TextValue "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "35.90"
Set-TextValue "E26" "  +7.34%  "
Set-TextValue "B27" "Filecoin"
Set-TextValue "C27" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D27" "10.05"
Set-TextValue "E27" "  -0.62%  "
Set-TextValue "E28" "  -1.57%  "
Set-TextValue "D29" "4.96"
Set-TextValue "E29" "  +4.20%  "
Set-TextValue "B30" "Toncoin"
Set-TextValue "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.79"
Set-TextValue "E30" "  +4.61%  "
Set-TextValue "B31" "Cosmos"
Set-TextValue "C31" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "12.26"
Set-TextValue "E31" "  -1.36%  "
Set-TextValue "B32" "Hedera"
Set-TextValue "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.119"
Set-TextValue "E32" "  +6.33%  "
Set-TextValue "D33" "7.13"
Set-TextValue "E33" "  -5.11%  "
Set-TextValue "B34" "InjectiveProtocol"
Set-TextValue "C34" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D34" "40.45"
Set-TextValue "E34" "  +1.37%  "
Set-TextValue "B35" "Kaspa"
Set-TextValue "C35" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D35" "0.157"
Set-TextValue "E35" "  -7.00%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.11%  "
Set-TextValue "D37" "56.18"
Set-TextValue "D38" "0.0₃0797"
Set-TextValue "E38" "  +21.86%  "
Set-TextValue "E39" "  +1.75%  "
Set-TextValue "D40" "0.151"
Set-TextValue "E40" "  +10.39%  "
Set-TextValue "E41" "  -0.24%  "
Set-TextValue "D42" "149.32"
Set-TextValue "E42" "  +2.07%  "
Set-TextValue "E43" "  -3.13%  "
Set-TextValue "D44" "2.92"
Set-TextValue "E44" "  -5.54%  "
Set-TextValue "B45" "LidoDAOToken"
Set-TextValue "C45" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D45" "3.26"
Set-TextValue "E45" "  -1.61%  "
Set-TextValue "B46" "ThetaToken"
Set-TextValue "C46" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D46" "2.63"
Set-TextValue "E46" "  +13.24%  "
Set-TextValue "E47" "  +22.73%  "
Set-TextValue "D48" "4.26"
Set-TextValue "E48" "  -5.40%  "
Set-TextValue "E49" "  -4.21%  "
Set-TextValue "E50" "  -2.01%  "
Set-TextValue "E51" "  +13.51%  "
